$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.468704581260681
$ws.Range("B1").Value = 1.534737229347229
$ws.Range("C1").Value = 1.458185434341431
$ws.Range("D1").Value = 1.436941146850586
$ws.Range("E1").Value = 1.028352379798889
